$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark all "bg", "bi" and "dg" files as translated (Nao -> SIM)
$ws.Range("C9:C225").Value = "SIM"
$ws.Range("C231:C308").Value = "SIM"
$ws.Range("C653:C663").Value = "SIM"
